# Update the "想去人数" (F column) figures for the affected con/exhibition
# rows on both the "展览" and "全部类型" worksheets, reflecting refreshed
# counts from the latest data pull (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 45
    4  = 541
    9  = 4486
    10 = 4367
    11 = 8
    12 = 16
    13 = 139
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
